$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; all existing rows 12..128 shift down to 13..129
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new data record
$ws.Range("A12").Value = 6
$ws.Range("B12").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value2 = 44537
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100112029
$ws.Range("G12").Value = "Orégano"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 34
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = 8441
$ws.Range("N12").Value = "`$/docena de atados"
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 2814
$ws.Range("Q12").Value = 3
$ws.Range("R12").Value = "Hortaliza"
